$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.1666666666666667
$ws.Range("C2").Value2 = 0.6666666666666666
$ws.Range("P2").Value2 = 0.1666666666666667

$ws.Range("C3").Value2 = 0.1111111111111111
$ws.Range("P3").Value2 = 0.7777777777777778
$ws.Range("S3").Value2 = 0.1111111111111111

$ws.Range("J4").Value2 = 0.25
$ws.Range("P4").Value2 = 0.75

$ws.Range("J6").Value2 = 0.1
$ws.Range("Q6").Value2 = 0.3
$ws.Range("R6").Value2 = 0.2

$ws.Range("B7").Value2 = 0.3333333333333333
$ws.Range("J7").Value2 = 0.1666666666666667
$ws.Range("Q7").Value2 = 0.1666666666666667
$ws.Range("R7").Value2 = 0.3333333333333333

$ws.Range("F8").Value2 = 0.03225806451612903
$ws.Range("J8").Value2 = 0.09677419354838709
$ws.Range("Q8").Value2 = 0.4193548387096774
$ws.Range("R8").Value2 = 0.1612903225806452
$ws.Range("S8").Value2 = 0.2903225806451613

$ws.Range("F9").Value2 = 0.07692307692307693
$ws.Range("J9").Value2 = 0.1538461538461539
$ws.Range("Q9").Value2 = 0.07692307692307693
$ws.Range("R9").Value2 = 0.3076923076923077
$ws.Range("S9").Value2 = 0.3846153846153846

$ws.Range("B10").Value2 = 0.06779661016949153
$ws.Range("D10").Value2 = 0.03389830508474576
$ws.Range("F10").Value2 = 0.05084745762711865
$ws.Range("J10").Value2 = 0.1610169491525424
$ws.Range("Q10").Value2 = 0.2711864406779661
$ws.Range("R10").Value2 = 0.1101694915254237
$ws.Range("S10").Value2 = 0.3050847457627119

$ws.Range("G11").Value2 = 0.25
$ws.Range("K11").Value2 = 0.25
$ws.Range("L11").Value2 = 0.5

$ws.Range("G12").Value2 = 1

$ws.Range("H15").Value2 = 0.1176470588235294
$ws.Range("I15").Value2 = 0.05882352941176471
$ws.Range("J15").Value2 = 0.6470588235294118
$ws.Range("M15").Value2 = 0.05882352941176471
$ws.Range("O15").Value2 = 0.05882352941176471
$ws.Range("S15").Value2 = 0.05882352941176471

$ws.Range("I16").Value2 = 0.1666666666666667
$ws.Range("J16").Value2 = 0.5833333333333334
$ws.Range("O16").Value2 = 0.08333333333333333
$ws.Range("S16").Value2 = 0.1666666666666667

$ws.Range("H17").Value2 = 0.2
$ws.Range("I17").Value2 = 0.04
$ws.Range("J17").Value2 = 0.58
$ws.Range("K17").Value2 = 0.02
$ws.Range("O17").Value2 = 0.08
$ws.Range("S17").Value2 = 0.08

$ws.Range("H18").Value2 = 0.16
$ws.Range("I18").Value2 = 0.12
$ws.Range("J18").Value2 = 0.52
$ws.Range("K18").Value2 = 0.04
$ws.Range("O18").Value2 = 0.08
$ws.Range("S18").Value2 = 0.08

$ws.Range("F19").Value2 = 0.01538461538461539
$ws.Range("H19").Value2 = 0.2307692307692308
$ws.Range("I19").Value2 = 0.07692307692307693
$ws.Range("J19").Value2 = 0.5076923076923077
$ws.Range("K19").Value2 = 0.01538461538461539
$ws.Range("M19").Value2 = 0.03076923076923077
$ws.Range("O19").Value2 = 0.09230769230769231
$ws.Range("S19").Value2 = 0.03076923076923077
